# Applies scheduled market-price refresh to Pandaemonium_Profits sheets
# (commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1625.9231
$ws.Range("J112").Value = 1857.909
$ws.Range("L112").Value = 5573.727000000001
$ws.Range("N112").Value = -7789.727000000001

$ws.Range("H127").Value = 928.41174
$ws.Range("I127").Value = 549.8333
$ws.Range("J127").Value = 1009.5357
$ws.Range("K127").Value = 1649.4999
$ws.Range("L127").Value = 3028.6071
$ws.Range("M127").Value = 3310.5001
$ws.Range("N127").Value = -12948.6071

$ws.Range("H137").Value = 2185.04
$ws.Range("I137").Value = 1681.25
$ws.Range("J137").Value = 4200.2
$ws.Range("K137").Value = 5043.75
$ws.Range("L137").Value = 12600.6
$ws.Range("M137").Value = -2493.75
$ws.Range("N137").Value = -17700.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4186.5776
$ws.Range("I61").Value = 3255.9033
$ws.Range("J61").Value = 6247.357
$ws.Range("K61").Value = 3255.9033
$ws.Range("L61").Value = 6247.357
$ws.Range("M61").Value = -3043.9033
$ws.Range("N61").Value = -6671.357

$ws.Range("H63").Value = 2815.1428
$ws.Range("I63").Value = 1633.3334
$ws.Range("J63").Value = 3701.5
$ws.Range("K63").Value = 1633.3334
$ws.Range("L63").Value = 3701.5
$ws.Range("M63").Value = -947.3334
$ws.Range("N63").Value = -5073.5

$ws.Range("H66").Value = 2815.1428
$ws.Range("I66").Value = 1633.3334
$ws.Range("J66").Value = 3701.5
$ws.Range("K66").Value = 8166.666999999999
$ws.Range("L66").Value = 18507.5
$ws.Range("M66").Value = -4734.666999999999
$ws.Range("N66").Value = -25371.5

$ws.Range("H74").Value = 4654.6313
$ws.Range("I74").Value = 2095.742
$ws.Range("J74").Value = 15986.857
$ws.Range("K74").Value = 2095.742
$ws.Range("L74").Value = 15986.857
$ws.Range("M74").Value = -1221.742
$ws.Range("N74").Value = -17734.857

$ws.Range("H77").Value = 4654.6313
$ws.Range("I77").Value = 2095.742
$ws.Range("J77").Value = 15986.857
$ws.Range("K77").Value = 10478.71
$ws.Range("L77").Value = 79934.285
$ws.Range("M77").Value = -6110.710000000001
$ws.Range("N77").Value = -88670.285

$ws.Range("H136").Value = 4186.5776
$ws.Range("I136").Value = 3255.9033
$ws.Range("J136").Value = 6247.357
$ws.Range("K136").Value = 9767.7099
$ws.Range("L136").Value = 18742.071
$ws.Range("M136").Value = -7217.7099
$ws.Range("N136").Value = -23842.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1716.1471
$ws.Range("I134").Value = 1582.8462
$ws.Range("J134").Value = 2149.375
$ws.Range("K134").Value = 4748.5386
$ws.Range("L134").Value = 6448.125
$ws.Range("M134").Value = -2213.5386
$ws.Range("N134").Value = -11518.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2578.25
$ws.Range("I31").Value = 1402.5
$ws.Range("J31").Value = 3326.4546
$ws.Range("K31").Value = 1402.5
$ws.Range("L31").Value = 3326.4546
$ws.Range("M31").Value = -1107.5
$ws.Range("N31").Value = -3916.4546

$ws.Range("H34").Value = 2578.25
$ws.Range("I34").Value = 1402.5
$ws.Range("J34").Value = 3326.4546
$ws.Range("K34").Value = 1402.5
$ws.Range("L34").Value = 3326.4546
$ws.Range("M34").Value = -1200.5
$ws.Range("N34").Value = -3730.4546

$ws.Range("H58").Value = 2116126.8
$ws.Range("I58").Value = 3788897
$ws.Range("J58").Value = 3153.6843
$ws.Range("K58").Value = 3788897
$ws.Range("L58").Value = 3153.6843
$ws.Range("M58").Value = -3788694
$ws.Range("N58").Value = -3559.6843

$ws.Range("H132").Value = 2343.2622
$ws.Range("I132").Value = 2119.7908
$ws.Range("J132").Value = 2877.111
$ws.Range("K132").Value = 6359.3724
$ws.Range("L132").Value = 8631.332999999999
$ws.Range("M132").Value = -3829.3724
$ws.Range("N132").Value = -13691.333

$ws.Range("H134").Value = 1954.5745
$ws.Range("I134").Value = 1838.5143
$ws.Range("K134").Value = 5515.5429
$ws.Range("M134").Value = -2980.5429

$ws.Range("H136").Value = 2116126.8
$ws.Range("I136").Value = 3788897
$ws.Range("J136").Value = 3153.6843
$ws.Range("K136").Value = 11366691
$ws.Range("L136").Value = 9461.052899999999
$ws.Range("M136").Value = -11364141
$ws.Range("N136").Value = -14561.0529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 700
$ws.Range("J35").Value = 700
$ws.Range("L35").Value = 2100
$ws.Range("N35").Value = -2676

$ws.Range("H113").Value = 678.93335
$ws.Range("I113").Value = 701.4286
$ws.Range("J113").Value = 600.2
$ws.Range("K113").Value = 2104.2858
$ws.Range("L113").Value = 1800.6
$ws.Range("M113").Value = 65.71420000000035
$ws.Range("N113").Value = -6140.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2749.875
$ws.Range("J126").Value = 2999.8333
$ws.Range("L126").Value = 8999.499899999999
$ws.Range("N126").Value = -13939.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 265.42105
$ws.Range("I55").Value = 268.8889
$ws.Range("J55").Value = 262.3
$ws.Range("K55").Value = 268.8889
$ws.Range("L55").Value = 262.3
$ws.Range("M55").Value = -95.88889999999998
$ws.Range("N55").Value = -608.3

$ws.Range("H132").Value = 3484.449
$ws.Range("I132").Value = 3066.743
$ws.Range("J132").Value = 4528.7144
$ws.Range("K132").Value = 9200.228999999999
$ws.Range("L132").Value = 13586.1432
$ws.Range("M132").Value = -6670.228999999999
$ws.Range("N132").Value = -18646.1432

$ws.Range("H136").Value = 3994.5
$ws.Range("I136").Value = 2398.9375
$ws.Range("J136").Value = 5958.269
$ws.Range("K136").Value = 7196.8125
$ws.Range("L136").Value = 17874.807
$ws.Range("M136").Value = -4646.8125
$ws.Range("N136").Value = -22974.807

$ws.Range("H140").Value = 60776.8
$ws.Range("J140").Value = 60776.8
$ws.Range("L140").Value = 60776.8
$ws.Range("N140").Value = -71136.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1684.674
$ws.Range("I132").Value = 857.23334
$ws.Range("J132").Value = 3236.125
$ws.Range("K132").Value = 2571.70002
$ws.Range("L132").Value = 9708.375
$ws.Range("M132").Value = -41.70002000000022
$ws.Range("N132").Value = -14768.375

$ws.Range("H136").Value = 5508.9614
$ws.Range("I136").Value = 2964.5881
$ws.Range("J136").Value = 10315
$ws.Range("K136").Value = 8893.764299999999
$ws.Range("L136").Value = 30945
$ws.Range("M136").Value = -6343.764299999999
$ws.Range("N136").Value = -36045
